$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Append a new row at the end of the table (after "errorDistributionID")
# and populate its cells to describe the "step" field.
$newRow = $t.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "step"
$newRow.Cells.Item(2).Range.Text = "Etape d'intégration du message"
$newRow.Cells.Item(3).Range.Text = "string"
$newRow.Cells.Item(4).Range.Text = "0..1"
$newRow.Cells.Item(5).Range.Text = "Nomenclature permettant d'identifier les différentes étapes d'intégration et de consultation du dossier dans le système émetteur"
